$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 560.2632
$ws.Range("I2").Value = 138.5
$ws.Range("J2").Value = 1028.8889
$ws.Range("K2").Value = 138.5
$ws.Range("L2").Value = 1028.8889
$ws.Range("M2").Value = -25.5
$ws.Range("N2").Value = -1254.8889
$ws.Range("H5").Value = 68.5
$ws.Range("I5").Value = 68.5
$ws.Range("K5").Value = 68.5
$ws.Range("M5").Value = 46.5
$ws.Range("H9").Value = 9463.583000000001
$ws.Range("I9").Value = 11189.223
$ws.Range("K9").Value = 11189.223
$ws.Range("M9").Value = -11020.223
$ws.Range("H32").Value = 17211.125
$ws.Range("I32").Value = 9888
$ws.Range("K32").Value = 9888
$ws.Range("M32").Value = -9562
$ws.Range("H43").Value = 10619.75
$ws.Range("I43").Value = 5999.5
$ws.Range("K43").Value = 5999.5
$ws.Range("M43").Value = -5930.5
$ws.Range("H55").Value = 42.142857
$ws.Range("I55").Value = 55
$ws.Range("K55").Value = 55
$ws.Range("M55").Value = 159
$ws.Range("H58").Value = 4243.5713
$ws.Range("J58").Value = 9611.5
$ws.Range("L58").Value = 28834.5
$ws.Range("N58").Value = -29134.5
$ws.Range("H64").Value = 4782.8335
$ws.Range("I64").Value = 4348.5
$ws.Range("K64").Value = 4348.5
$ws.Range("M64").Value = -4100.5
$ws.Range("H67").Value = 4782.8335
$ws.Range("I67").Value = 4348.5
$ws.Range("K67").Value = 4348.5
$ws.Range("M67").Value = -3490.5
$ws.Range("H70").Value = 4990.9287
$ws.Range("I70").Value = 2953.2
$ws.Range("J70").Value = 7342.154
$ws.Range("K70").Value = 8859.599999999999
$ws.Range("L70").Value = 22026.462
$ws.Range("M70").Value = -8589.599999999999
$ws.Range("N70").Value = -22566.462
$ws.Range("H73").Value = 4990.9287
$ws.Range("I73").Value = 2953.2
$ws.Range("J73").Value = 7342.154
$ws.Range("K73").Value = 8859.599999999999
$ws.Range("L73").Value = 22026.462
$ws.Range("M73").Value = -7923.599999999999
$ws.Range("N73").Value = -23898.462
$ws.Range("H106").Value = 3223.0667
$ws.Range("I106").Value = 2490.4
$ws.Range("J106").Value = 4688.4
$ws.Range("K106").Value = 2490.4
$ws.Range("L106").Value = 4688.4
$ws.Range("M106").Value = -1859.4
$ws.Range("N106").Value = -5950.4
$ws.Range("H112").Value = 43059.38
$ws.Range("I112").Value = 85217.164
$ws.Range("J112").Value = 29746.395
$ws.Range("K112").Value = 255651.492
$ws.Range("L112").Value = 89239.185
$ws.Range("M112").Value = -254543.492
$ws.Range("N112").Value = -91455.185
$ws.Range("H135").Value = 115385450
$ws.Range("I135").Value = 45455120
$ws.Range("K135").Value = 409096080
$ws.Range("M135").Value = -409093545
$ws.Range("H137").Value = 1761.5
$ws.Range("I137").Value = 1081.75
$ws.Range("K137").Value = 3245.25
$ws.Range("M137").Value = -695.25
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 670008
$ws.Range("I2").Value = 818205.2
$ws.Range("J2").Value = 3120.75
$ws.Range("K2").Value = 818205.2
$ws.Range("L2").Value = 3120.75
$ws.Range("M2").Value = -818092.2
$ws.Range("N2").Value = -3346.75
$ws.Range("H62").Value = 24666
$ws.Range("J62").Value = 24666
$ws.Range("L62").Value = 24666
$ws.Range("N62").Value = -25914
$ws.Range("H65").Value = 24666
$ws.Range("J65").Value = 24666
$ws.Range("L65").Value = 73998
$ws.Range("N65").Value = -80238
$ws.Range("H116").Value = 670008
$ws.Range("I116").Value = 818205.2
$ws.Range("J116").Value = 3120.75
$ws.Range("K116").Value = 818205.2
$ws.Range("L116").Value = 3120.75
$ws.Range("M116").Value = -815911.2
$ws.Range("N116").Value = -7708.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 670008
$ws.Range("I3").Value = 818205.2
$ws.Range("J3").Value = 3120.75
$ws.Range("K3").Value = 818205.2
$ws.Range("L3").Value = 3120.75
$ws.Range("M3").Value = -818091.2
$ws.Range("N3").Value = -3348.75
$ws.Range("H16").Value = 608
$ws.Range("I16").Value = 608
$ws.Range("K16").Value = 608
$ws.Range("M16").Value = -438
$ws.Range("H94").Value = 915.1667
$ws.Range("I94").Value = 886.4375
$ws.Range("J94").Value = 1145
$ws.Range("K94").Value = 886.4375
$ws.Range("L94").Value = 1145
$ws.Range("M94").Value = -435.4375
$ws.Range("N94").Value = -2047
$ws.Range("H105").Value = 1615.9584
$ws.Range("I105").Value = 1562.0952
$ws.Range("J105").Value = 1993
$ws.Range("K105").Value = 1562.0952
$ws.Range("L105").Value = 1993
$ws.Range("M105").Value = 184.9048
$ws.Range("N105").Value = -5487
$ws.Range("H128").Value = 4999
$ws.Range("I128").Value = 4999
$ws.Range("K128").Value = 14997
$ws.Range("M128").Value = -12507
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 169.3
$ws.Range("I7").Value = 71.5
$ws.Range("J7").Value = 267.1
$ws.Range("K7").Value = 71.5
$ws.Range("L7").Value = 267.1
$ws.Range("M7").Value = 41.5
$ws.Range("N7").Value = -493.1
$ws.Range("H107").Value = 915633.6
$ws.Range("I107").Value = 1208304.2
$ws.Range("K107").Value = 1208304.2
$ws.Range("M107").Value = -1206384.2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 133.5
$ws.Range("I2").Value = 124.916664
$ws.Range("K2").Value = 749.499984
$ws.Range("M2").Value = -636.499984
$ws.Range("H3").Value = 33000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 33000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 99000
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -99224
$ws.Range("H8").Value = 758
$ws.Range("I8").Value = 758
$ws.Range("K8").Value = 2274
$ws.Range("M8").Value = -2135
$ws.Range("H86").Value = 788
$ws.Range("I86").Value = 690
$ws.Range("K86").Value = 2070
$ws.Range("M86").Value = -884
$ws.Range("H87").Value = 19249.75
$ws.Range("I87").Value = 19249.75
$ws.Range("K87").Value = 57749.25
$ws.Range("M87").Value = -56501.25
$ws.Range("H89").Value = 788
$ws.Range("I89").Value = 690
$ws.Range("K89").Value = 6210
$ws.Range("M89").Value = -282
$ws.Range("H90").Value = 19249.75
$ws.Range("I90").Value = 19249.75
$ws.Range("K90").Value = 173247.75
$ws.Range("M90").Value = -167007.75
$ws.Range("H123").Value = 5301.875
$ws.Range("I123").Value = 924
$ws.Range("J123").Value = 6761.1665
$ws.Range("K123").Value = 2772
$ws.Range("L123").Value = 20283.4995
$ws.Range("M123").Value = -322
$ws.Range("N123").Value = -25183.4995
$ws.Range("H126").Value = 1264.5
$ws.Range("I126").Value = 1264.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3793.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 1146.5
$ws.Range("N126").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4999
$ws.Range("I80").Value = 4998.5
$ws.Range("K80").Value = 4998.5
$ws.Range("M80").Value = -4000.5
$ws.Range("H83").Value = 4999
$ws.Range("I83").Value = 4998.5
$ws.Range("K83").Value = 24992.5
$ws.Range("M83").Value = -20000.5
$ws.Range("H97").Value = 613.24
$ws.Range("I97").Value = 400.875
$ws.Range("J97").Value = 990.7778
$ws.Range("K97").Value = 400.875
$ws.Range("L97").Value = 990.7778
$ws.Range("M97").Value = 95.125
$ws.Range("N97").Value = -1982.7778
$ws.Range("H102").Value = 1735.8182
$ws.Range("I102").Value = 1545.4445
$ws.Range("J102").Value = 2592.5
$ws.Range("K102").Value = 1545.4445
$ws.Range("L102").Value = 2592.5
$ws.Range("M102").Value = 76.55549999999994
$ws.Range("N102").Value = -5836.5
$ws.Range("H122").Value = 178671.14
$ws.Range("I122").Value = 602849
$ws.Range("K122").Value = 1808547
$ws.Range("M122").Value = -1806097
$ws.Range("H139").Value = 127497.5
$ws.Range("J139").Value = 127497.5
$ws.Range("L139").Value = 127497.5
$ws.Range("N139").Value = -137777.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2078.1428
$ws.Range("I22").Value = 2109.6
$ws.Range("K22").Value = 2109.6
$ws.Range("M22").Value = -1814.6
$ws.Range("H27").Value = 2078.1428
$ws.Range("I27").Value = 2109.6
$ws.Range("K27").Value = 2109.6
$ws.Range("M27").Value = -2002.6
$ws.Range("H46").Value = 1517.8667
$ws.Range("I46").Value = 1666.5
$ws.Range("J46").Value = 1418.7778
$ws.Range("K46").Value = 1666.5
$ws.Range("L46").Value = 1418.7778
$ws.Range("M46").Value = -1478.5
$ws.Range("N46").Value = -1794.7778
$ws.Range("H47").Value = 25750
$ws.Range("J47").Value = 25750
$ws.Range("L47").Value = 25750
$ws.Range("N47").Value = -26730
$ws.Range("H52").Value = 25750
$ws.Range("J52").Value = 25750
$ws.Range("L52").Value = 25750
$ws.Range("N52").Value = -26216
$ws.Range("H68").Value = 3572500
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 3572500
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H107").Value = 3499.5
$ws.Range("I107").Value = 3499.5
$ws.Range("K107").Value = 3499.5
$ws.Range("M107").Value = -1579.5
$ws.Range("H140").Value = 64933.727
$ws.Range("J140").Value = 64933.727
$ws.Range("L140").Value = 64933.727
$ws.Range("N140").Value = -75293.727
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 42883884
$ws.Range("J4").Value = 60022440
$ws.Range("L4").Value = 60022440
$ws.Range("N4").Value = -60022666
